$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Insert three new sheets (weather, report, animals) immediately before the
# existing "customers" sheet, preserving creation order so sheetId values
# come out as weather=8, report=9, animals=10 (matches target workbook.xml).
# ---------------------------------------------------------------------------
$new_students = $wb.Worksheets.Item("new_students")

$weather = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $new_students)
$weather.Name = "weather"

$report = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weather)
$report.Name = "report"

$animals = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $report)
$animals.Name = "animals"

# ---------------------------------------------------------------------------
# weather sheet: city / month / temperature
# ---------------------------------------------------------------------------
$weather.Cells.Item(1, 1).Value = "city"
$weather.Cells.Item(1, 2).Value = "month"
$weather.Cells.Item(1, 3).Value = "temperature"

$weatherData = @(
  @("Jacksonville", "January", 13),
  @("Jacksonville", "February", 23),
  @("Jacksonville", "March", 38),
  @("Jacksonville", "April", 5),
  @("Jacksonville", "May", 34),
  @("ElPaso", "January", 20),
  @("ElPaso", "February", 6),
  @("ElPaso", "March", 26),
  @("ElPaso", "April", 2),
  @("ElPaso", "May", 43)
)

$r = 2
foreach ($row in $weatherData) {
  $weather.Cells.Item($r, 1).Value = $row[0]
  $weather.Cells.Item($r, 2).Value = $row[1]
  $weather.Cells.Item($r, 3).Value = $row[2]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# report sheet: product / quarter_1 / quarter_2 / quarter_3 / quarter_4
# ---------------------------------------------------------------------------
$report.Cells.Item(1, 1).Value = "product"
$report.Cells.Item(1, 2).Value = "quarter_1"
$report.Cells.Item(1, 3).Value = "quarter_2"
$report.Cells.Item(1, 4).Value = "quarter_3"
$report.Cells.Item(1, 5).Value = "quarter_4"

$reportData = @(
  @("Umbrella", 417, 224, 379, 611),
  @("SleepingBag", 800, 936, 93, 875)
)

$r = 2
foreach ($row in $reportData) {
  $report.Cells.Item($r, 1).Value = $row[0]
  $report.Cells.Item($r, 2).Value = $row[1]
  $report.Cells.Item($r, 3).Value = $row[2]
  $report.Cells.Item($r, 4).Value = $row[3]
  $report.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# animals sheet: name / species / age / weight
# ---------------------------------------------------------------------------
$animals.Cells.Item(1, 1).Value = "name"
$animals.Cells.Item(1, 2).Value = "species"
$animals.Cells.Item(1, 3).Value = "age"
$animals.Cells.Item(1, 4).Value = "weight"

$animalsData = @(
  @("Tatiana", "Snake", 98, 464),
  @("Khaled", "Giraffe", 50, 41),
  @("Alex", "Leopard", 6, 328),
  @("Jonathan", "Monkey", 45, 463),
  @("Stefan", "Bear", 100, 50),
  @("Tommy", "Panda", 26, 349)
)

$r = 2
foreach ($row in $animalsData) {
  $animals.Cells.Item($r, 1).Value = $row[0]
  $animals.Cells.Item($r, 2).Value = $row[1]
  $animals.Cells.Item($r, 3).Value = $row[2]
  $animals.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# Selections for the touched sheets.
# ---------------------------------------------------------------------------
$new_students.Activate()
[void]$new_students.Range("E18").Select()

$weather.Activate()
[void]$weather.Range("E18").Select()

$report.Activate()
[void]$report.Range("B3").Select()

$animals.Activate()
[void]$animals.Range("H20").Select()
